$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF") should use the same bold/
# bordered/centered header style as the existing H1 ("IP") header.
# Copying H1's format onto I1:J1 reuses the existing style entry instead
# of minting new ones, then the text values are written on top.
$ws.Range("H1").Copy($ws.Range("I1:J1"))
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows 2-11: columns I (I0) and J (IF) hold the same numeric value.
$ws.Range("I2").Value = 8
$ws.Range("J2").Value = 8

$ws.Range("I3").Value = 5
$ws.Range("J3").Value = 5

$ws.Range("I4").Value = 5
$ws.Range("J4").Value = 5

$ws.Range("I5").Value = 10
$ws.Range("J5").Value = 10

$ws.Range("I6").Value = 8
$ws.Range("J6").Value = 8

$ws.Range("I7").Value = 4
$ws.Range("J7").Value = 4

$ws.Range("I8").Value = 5
$ws.Range("J8").Value = 5

$ws.Range("I9").Value = 4
$ws.Range("J9").Value = 4

$ws.Range("I10").Value = 4
$ws.Range("J10").Value = 4

$ws.Range("I11").Value = 8
$ws.Range("J11").Value = 8
